# Add NSW "Anzac Day (additional)" observance rows for 2026 and 2027
# (Anzac Day falls on a weekend in those years, so NSW observes an
# additional holiday on the following Monday).
#
# Insert from the bottom of the sheet upward so earlier row numbers
# stay valid while we work.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 2027: Anzac Day (VIC) is on row 43 -> insert new row 44 ---
$ws.Rows.Item(44).Insert()
# Force the date-like text to stay literal text (not auto-converted to a
# date serial) by applying a Text number format before assigning it, then
# strip the formatting again so the new cell matches its plain neighbours.
$ws.Cells.Item(44, 1).NumberFormat = "@"
$ws.Cells.Item(44, 1).Value = "2027-04-26"
$ws.Cells.Item(44, 1).ClearFormats()
$ws.Cells.Item(44, 2).Value = "NSW"
$ws.Cells.Item(44, 3).Value = "Anzac Day (additional)"

# --- 2026: Anzac Day (VIC) is on row 16 -> insert new row 17 ---
$ws.Rows.Item(17).Insert()
$ws.Cells.Item(17, 1).NumberFormat = "@"
$ws.Cells.Item(17, 1).Value = "2026-04-27"
$ws.Cells.Item(17, 1).ClearFormats()
$ws.Cells.Item(17, 2).Value = "NSW"
$ws.Cells.Item(17, 3).Value = "Anzac Day (additional)"
